$d = $word.ActiveDocument

$replacements = @(
    @{old = "647÷4=161, 3"; new = "255÷5=51, 0"},
    @{old = "635÷4=158, 3"; new = "848÷5=169, 3"},
    @{old = "358÷9=39, 7"; new = "702÷6=117, 0"},
    @{old = "565÷4=141, 1"; new = "916÷8=114, 4"},
    @{old = "578÷3=192, 2"; new = "422÷3=140, 2"},
    @{old = "227÷6=37, 5"; new = "242÷7=34, 4"},
    @{old = "857÷3=285, 2"; new = "400÷3=133, 1"},
    @{old = "539÷6=89, 5"; new = "155÷7=22, 1"},
    @{old = "384÷9=42, 6"; new = "971÷4=242, 3"},
    @{old = "432÷4=108, 0"; new = "135÷7=19, 2"},
    @{old = "751÷7=107, 2"; new = "928÷6=154, 4"},
    @{old = "677÷4=169, 1"; new = "101÷6=16, 5"},
    @{old = "152÷4=38, 0"; new = "656÷4=164, 0"},
    @{old = "695÷3=231, 2"; new = "157÷8=19, 5"},
    @{old = "137÷8=17, 1"; new = "759÷6=126, 3"},
    @{old = "591÷5=118, 1"; new = "877÷7=125, 2"},
    @{old = "127÷2=63, 1"; new = "153÷2=76, 1"},
    @{old = "625÷9=69, 4"; new = "937÷9=104, 1"},
    @{old = "146÷6=24, 2"; new = "150÷6=25, 0"},
    @{old = "834÷4=208, 2"; new = "962÷7=137, 3"},
    @{old = "404÷2=202, 0"; new = "702÷4=175, 2"},
    @{old = "657÷7=93, 6"; new = "578÷5=115, 3"},
    @{old = "421÷8=52, 5"; new = "270÷7=38, 4"},
    @{old = "623÷4=155, 3"; new = "698÷3=232, 2"},
    @{old = "178÷7=25, 3"; new = "299÷3=99, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
